$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 781.3333
$ws.Range("I34").Value = 172
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 172
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = 31
$ws.Range("N34").Value = -2406
# Row 36
$ws.Range("H36").Value = 781.3333
$ws.Range("I36").Value = 172
$ws.Range("J36").Value = 2000
$ws.Range("K36").Value = 172
$ws.Range("L36").Value = 2000
$ws.Range("M36").Value = 543
$ws.Range("N36").Value = -3430
# Row 80
$ws.Range("H80").Value = 600.625
$ws.Range("I80").Value = 525
$ws.Range("K80").Value = 1575
$ws.Range("M80").Value = -577
# Row 83
$ws.Range("H83").Value = 600.625
$ws.Range("I83").Value = 525
$ws.Range("K83").Value = 4725
$ws.Range("M83").Value = 267
# Row 116
$ws.Range("H116").Value = 4003
$ws.Range("J116").Value = 3006
$ws.Range("L116").Value = 3006
$ws.Range("N116").Value = -9890
# Row 125
$ws.Range("H125").Value = 666.5
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920
# Row 132
$ws.Range("H132").Value = 8228.833000000001
$ws.Range("I132").Value = 3936.5
$ws.Range("J132").Value = 10375
$ws.Range("K132").Value = 11809.5
$ws.Range("L132").Value = 31125
$ws.Range("M132").Value = -9279.5
$ws.Range("N132").Value = -36185
# Row 138
$ws.Range("H138").Value = 6571.2856
$ws.Range("J138").Value = 4600
$ws.Range("L138").Value = 13800
$ws.Range("N138").Value = -24080
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 829.5454999999999
$ws.Range("I2").Value = 952.7778
$ws.Range("J2").Value = 275
$ws.Range("K2").Value = 952.7778
$ws.Range("L2").Value = 275
$ws.Range("M2").Value = -839.7778
$ws.Range("N2").Value = -501
# Row 32
$ws.Range("H32").Value = 6907.032
$ws.Range("I32").Value = 4968.074
$ws.Range("K32").Value = 4968.074
$ws.Range("M32").Value = -4681.074
# Row 61
$ws.Range("H61").Value = 1997.5
$ws.Range("I61").Value = 1997
$ws.Range("J61").Value = 1998
$ws.Range("K61").Value = 1997
$ws.Range("L61").Value = 1998
$ws.Range("M61").Value = -1785
$ws.Range("N61").Value = -2422
# Row 74
$ws.Range("H74").Value = 1967.8
$ws.Range("I74").Value = 1984.7778
$ws.Range("K74").Value = 1984.7778
$ws.Range("M74").Value = -1110.7778
# Row 77
$ws.Range("H77").Value = 1967.8
$ws.Range("I77").Value = 1984.7778
$ws.Range("K77").Value = 9923.889000000001
$ws.Range("M77").Value = -5555.889000000001
# Row 102
$ws.Range("H102").Value = 1563.3334
$ws.Range("I102").Value = 1945
$ws.Range("K102").Value = 1945
$ws.Range("M102").Value = -323
# Row 116
$ws.Range("H116").Value = 829.5454999999999
$ws.Range("I116").Value = 952.7778
$ws.Range("J116").Value = 275
$ws.Range("K116").Value = 952.7778
$ws.Range("L116").Value = 275
$ws.Range("M116").Value = 1341.2222
$ws.Range("N116").Value = -4863
# Row 122
$ws.Range("H122").Value = 3324.6667
$ws.Range("I122").Value = 2990.25
$ws.Range("K122").Value = 8970.75
$ws.Range("M122").Value = -6520.75
# Row 132
$ws.Range("H132").Value = 1537.125
$ws.Range("I132").Value = 1474.75
$ws.Range("K132").Value = 4424.25
$ws.Range("M132").Value = -1894.25
# Row 136
$ws.Range("H136").Value = 1997.5
$ws.Range("I136").Value = 1997
$ws.Range("J136").Value = 1998
$ws.Range("K136").Value = 5991
$ws.Range("L136").Value = 5994
$ws.Range("M136").Value = -3441
$ws.Range("N136").Value = -11094

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 829.5454999999999
$ws.Range("I3").Value = 952.7778
$ws.Range("J3").Value = 275
$ws.Range("K3").Value = 952.7778
$ws.Range("L3").Value = 275
$ws.Range("M3").Value = -838.7778
$ws.Range("N3").Value = -503
# Row 20
$ws.Range("H20").Value = 3308.3333
$ws.Range("I20").Value = 3308.3333
$ws.Range("K20").Value = 3308.3333
$ws.Range("M20").Value = -3061.3333
# Row 105
$ws.Range("H105").Value = 2219.8333
$ws.Range("I105").Value = 1829.75
$ws.Range("K105").Value = 1829.75
$ws.Range("M105").Value = -82.75
# Row 107
$ws.Range("H107").Value = 5683.4
$ws.Range("J107").Value = 1400
$ws.Range("L107").Value = 1400
$ws.Range("N107").Value = -5240
# Row 134
$ws.Range("H134").Value = 6498.857
$ws.Range("I134").Value = 3198.7
$ws.Range("K134").Value = 9596.099999999999
$ws.Range("M134").Value = -7061.099999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 118.5
$ws.Range("J2").Value = 233
$ws.Range("L2").Value = 233
$ws.Range("N2").Value = -459
# Row 31
$ws.Range("H31").Value = 18421.875
$ws.Range("I31").Value = 18925
$ws.Range("J31").Value = 14900
$ws.Range("K31").Value = 18925
$ws.Range("L31").Value = 14900
$ws.Range("M31").Value = -18630
$ws.Range("N31").Value = -15490
# Row 34
$ws.Range("H34").Value = 18421.875
$ws.Range("I34").Value = 18925
$ws.Range("J34").Value = 14900
$ws.Range("K34").Value = 18925
$ws.Range("L34").Value = 14900
$ws.Range("M34").Value = -18723
$ws.Range("N34").Value = -15304
# Row 58
$ws.Range("H58").Value = 2545.6667
$ws.Range("I58").Value = 2254.8
$ws.Range("J58").Value = 4000
$ws.Range("K58").Value = 2254.8
$ws.Range("L58").Value = 4000
$ws.Range("M58").Value = -2051.8
$ws.Range("N58").Value = -4406
# Row 86
$ws.Range("H86").Value = 2869.2856
$ws.Range("J86").Value = 2999
$ws.Range("L86").Value = 2999
$ws.Range("N86").Value = -5245
# Row 89
$ws.Range("H89").Value = 2869.2856
$ws.Range("J89").Value = 2999
$ws.Range("L89").Value = 14995
$ws.Range("N89").Value = -26227
# Row 136
$ws.Range("H136").Value = 2545.6667
$ws.Range("I136").Value = 2254.8
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6764.400000000001
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4214.400000000001
$ws.Range("N136").Value = -17100
# Row 141
$ws.Range("H141").Value = 408137.72
$ws.Range("I141").Value = 88796
$ws.Range("J141").Value = 461361.34
$ws.Range("K141").Value = 88796
$ws.Range("L141").Value = 461361.34
$ws.Range("M141").Value = -83616
$ws.Range("N141").Value = -471721.34

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 4624.6665
$ws.Range("I122").Value = 4743.5
$ws.Range("J122").Value = 4529.6
$ws.Range("K122").Value = 42691.5
$ws.Range("L122").Value = 40766.4
$ws.Range("M122").Value = -40241.5
$ws.Range("N122").Value = -45666.4
# Row 129
$ws.Range("H129").Value = 3387.375
$ws.Range("J129").Value = 250
$ws.Range("L129").Value = 750
$ws.Range("N129").Value = -10750
# Row 132
$ws.Range("H132").Value = 4000
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 36000
$ws.Range("M132").Value = -33470
# Row 140
$ws.Range("H140").Value = 2103.6667
$ws.Range("I140").Value = 2103.6667
$ws.Range("K140").Value = 6311.000100000001
$ws.Range("M140").Value = -1131.000100000001
# Row 141
$ws.Range("H141").Value = 2379.8333
$ws.Range("I141").Value = 2379.8333
$ws.Range("K141").Value = 7139.499899999999
$ws.Range("M141").Value = -1959.499899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 9997.5
$ws.Range("J102").Value = 9995
$ws.Range("L102").Value = 9995
$ws.Range("N102").Value = -13239
# Row 107
$ws.Range("H107").Value = 820.4
$ws.Range("I107").Value = 700.5
$ws.Range("K107").Value = 700.5
$ws.Range("M107").Value = 1219.5
# Row 113
$ws.Range("H113").Value = 945
$ws.Range("I113").Value = 860
$ws.Range("K113").Value = 860
$ws.Range("M113").Value = 1310
# Row 122
$ws.Range("H122").Value = 11415.765
$ws.Range("I122").Value = 7226.5835
$ws.Range("J122").Value = 21469.8
$ws.Range("K122").Value = 21679.7505
$ws.Range("L122").Value = 64409.39999999999
$ws.Range("M122").Value = -19229.7505
$ws.Range("N122").Value = -69309.39999999999
# Row 132
$ws.Range("H132").Value = 2368.5715
$ws.Range("I132").Value = 2096.75
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 6290.25
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -3760.25
$ws.Range("N132").Value = -17058.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("I16").Value = 325
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 325
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -155
$ws.Range("N16").ClearContents()
# Row 132
$ws.Range("H132").Value = 24455.334
$ws.Range("I132").Value = 27599.834
$ws.Range("K132").Value = 82799.50199999999
$ws.Range("M132").Value = -80269.50199999999
# Row 136
$ws.Range("H136").Value = 7150.1665
$ws.Range("I136").Value = 3875.25
$ws.Range("K136").Value = 11625.75
$ws.Range("M136").Value = -9075.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1160.6
$ws.Range("I126").Value = 1160.6
$ws.Range("K126").Value = 3481.8
$ws.Range("M126").Value = -1011.8
# Row 132
$ws.Range("H132").Value = 6539.3335
$ws.Range("I132").Value = 3240.923
$ws.Range("J132").Value = 11899.25
$ws.Range("K132").Value = 9722.769
$ws.Range("L132").Value = 35697.75
$ws.Range("M132").Value = -7192.769
$ws.Range("N132").Value = -40757.75
# Row 136
$ws.Range("H136").Value = 1780.1666
$ws.Range("I136").Value = 1780.1666
$ws.Range("K136").Value = 5340.4998
$ws.Range("M136").Value = -2790.4998
